$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold the restaurant name.
$ws.Columns("A").Insert()

# Header for the new column.
$ws.Range("A1").Value = "restaurant"

# Fill the new column with the restaurant name for every data row.
$ws.Range("A2:A180").Value = "Wendys"

# Match the author's final selection state (cosmetic, harmless).
$ws.Range("H11").Select() | Out-Null
